# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and fix the swapped WEMIXTOKEN / VeChain rows (38-39).

function Set-TextValue($cell, $value) {
    # Force the cell to be written as text (matches the sheet's existing
    # inline-string cells), even when the value looks numeric
    # (e.g. "43.31" or "1.061"), then drop the temporary number format
    # again so the cell keeps its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "21.766.04"
Set-TextValue $ws.Range("E2") "  -1.68%  "
Set-TextValue $ws.Range("D3") "1.540.83"
Set-TextValue $ws.Range("E3") "  -1.36%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("E5") "  +0.06%  "
Set-TextValue $ws.Range("D6") "290.23"
Set-TextValue $ws.Range("E6") "  +0.22%  "
Set-TextValue $ws.Range("D7") "0.3894"
Set-TextValue $ws.Range("E7") "  +2.47%  "
Set-TextValue $ws.Range("D8") "0.3192"
Set-TextValue $ws.Range("E8") "  -3.07%  "
Set-TextValue $ws.Range("D9") "43.31"
Set-TextValue $ws.Range("E9") "  -0.56%  "
Set-TextValue $ws.Range("D10") "0.07191"
Set-TextValue $ws.Range("E10") "  -2.56%  "
Set-TextValue $ws.Range("D11") "1.061"
Set-TextValue $ws.Range("E11") "  -7.15%  "
Set-TextValue $ws.Range("D12") "1.002"
Set-TextValue $ws.Range("E12") "  +0.09%  "
Set-TextValue $ws.Range("D13") "5.636"
Set-TextValue $ws.Range("E13") "  -3.42%  "
Set-TextValue $ws.Range("D14") "18.57"
Set-TextValue $ws.Range("E14") "  -7.37%  "
Set-TextValue $ws.Range("D15") "6.621"
Set-TextValue $ws.Range("E15") "  -4.01%  "
Set-TextValue $ws.Range("D16") "1.542.68"
Set-TextValue $ws.Range("E16") "  -1.09%  "
Set-TextValue $ws.Range("E17") "  +0.60%  "
Set-TextValue $ws.Range("D18") "0.06584"
Set-TextValue $ws.Range("E18") "  -1.16%  "
Set-TextValue $ws.Range("D19") "83.51"
Set-TextValue $ws.Range("E19") "  -2.70%  "
Set-TextValue $ws.Range("D20") "1.001"
Set-TextValue $ws.Range("E20") "  +0.08%  "
Set-TextValue $ws.Range("D21") "6.142"
Set-TextValue $ws.Range("E21") "  -5.12%  "
Set-TextValue $ws.Range("D22") "15.38"
Set-TextValue $ws.Range("E22") "  -5.02%  "
Set-TextValue $ws.Range("D23") "10.88"
Set-TextValue $ws.Range("E23") "  -7.31%  "
Set-TextValue $ws.Range("D24") "2.368"
Set-TextValue $ws.Range("E24") "  +4.38%  "
Set-TextValue $ws.Range("D25") "21.758.31"
Set-TextValue $ws.Range("E25") "  -1.78%  "
Set-TextValue $ws.Range("D26") "2.392"
Set-TextValue $ws.Range("E26") "  -6.65%  "
Set-TextValue $ws.Range("D27") "145.29"
Set-TextValue $ws.Range("E27") "  -4.07%  "
Set-TextValue $ws.Range("D28") "18.38"
Set-TextValue $ws.Range("E28") "  -4.04%  "
Set-TextValue $ws.Range("D29") "4.850"
Set-TextValue $ws.Range("E29") "  -0.16%  "
Set-TextValue $ws.Range("D30") "1.719.60"
Set-TextValue $ws.Range("E30") "  -0.96%  "
Set-TextValue $ws.Range("D31") "117.60"
Set-TextValue $ws.Range("E31") "  -3.21%  "
Set-TextValue $ws.Range("D32") "0.9630"
Set-TextValue $ws.Range("E32") "  -14.37%  "
Set-TextValue $ws.Range("D33") "5.874"
Set-TextValue $ws.Range("E33") "  -3.10%  "
Set-TextValue $ws.Range("D34") "0.08222"
Set-TextValue $ws.Range("E34") "  +0.49%  "
Set-TextValue $ws.Range("D35") "8.961"
Set-TextValue $ws.Range("E35") "  -4.66%  "
Set-TextValue $ws.Range("D36") "0.06103"
Set-TextValue $ws.Range("E36") "  -1.94%  "
Set-TextValue $ws.Range("D37") "5.134"
Set-TextValue $ws.Range("E37") "  -3.15%  "
Set-TextValue $ws.Range("B38") "WEMIXTOKEN"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D38") "1.481"
Set-TextValue $ws.Range("E38") "  -20.55%  "
Set-TextValue $ws.Range("B39") "VeChain"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.02205"
Set-TextValue $ws.Range("E39") "  -4.45%  "
Set-TextValue $ws.Range("D40") "0.2041"
Set-TextValue $ws.Range("E40") "  -4.81%  "
Set-TextValue $ws.Range("D41") "1.186"
Set-TextValue $ws.Range("E41") "  -4.35%  "
Set-TextValue $ws.Range("E42") "  +0.01%  "
Set-TextValue $ws.Range("E43") "  -3.84%  "
Set-TextValue $ws.Range("D44") "0.5744"
Set-TextValue $ws.Range("E44") "  -4.31%  "
Set-TextValue $ws.Range("D45") "13.08"
Set-TextValue $ws.Range("E45") "  -4.94%  "
Set-TextValue $ws.Range("D46") "3.740"
Set-TextValue $ws.Range("E46") "  -0.32%  "
Set-TextValue $ws.Range("D47") "0.5511"
Set-TextValue $ws.Range("E47") "  -5.06%  "
Set-TextValue $ws.Range("D48") "118.31"
Set-TextValue $ws.Range("E48") "  -2.27%  "
Set-TextValue $ws.Range("D49") "1.869"
Set-TextValue $ws.Range("E49") "  -5.72%  "
Set-TextValue $ws.Range("D50") "1.140"
Set-TextValue $ws.Range("E50") "  -2.81%  "
Set-TextValue $ws.Range("D51") "0.06733"
Set-TextValue $ws.Range("E51") "  -3.61%  "
